# "Realizado 6.1 'Situacion 1'"
# Fill in the "Todos" answer for the row that was left blank, and leave
# the cursor where the author ended up (cell E12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("D11")

# Give D11 the same "filled" look as the other answered rows in this
# column (D3, D4, D7, D9 all use the "40% - Énfasis5" cell style).
$target.Style = "40% - Énfasis5"

# Record the answer.
$target.Value = "Todos"

# Match the author's final selection/cursor position.
$ws.Range("E12").Select()
